# submission clearance scripts updates
# - Adds three new "Clearances" test-data blocks to the InsuredPageData sheet
#   (testClearancesSubmissionFunctionality, testCancelClearancesFunctionality,
#   testSubmissionClearanceComplete)
# - Updates the "testCheckDuplicateSubmission" sample row (C23/D23) to use the
#   new "venkat qa" / "venkatqa.com" test fixtures instead of the old
#   "Glasscock Chevrolet, Inc." / profrisk.com values
# - Makes InsuredPageData the active/selected sheet (was DashboardPageData)

$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item(1)
$wsDashboard = $wb.Worksheets.Item(2)
$wsInsured = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Update the existing duplicate-submission sample data (row 23)
# ---------------------------------------------------------------------------
$wsInsured.Range("C23").Value = "venkat qa"
$wsInsured.Range("D23").Value = "venkatqa.com"

# ---------------------------------------------------------------------------
# 2. Add the three new clearance test-data blocks, reusing the formatting of
#    the existing "testCheckDuplicateSubmission" block (rows 21-23) as a
#    template so borders/fill/etc. stay consistent with the rest of the sheet
# ---------------------------------------------------------------------------
$templateRange = $wsInsured.Range("A21:I23")

# --- testClearancesSubmissionFunctionality (rows 26-28) -------------------
$templateRange.Copy() | Out-Null
$wsInsured.Range("A26").PasteSpecial(-4122) | Out-Null

$wsInsured.Range("A26").Value = "testClearancesSubmissionFunctionality"

$wsInsured.Range("A27").Value = "runMode"
$wsInsured.Range("B27").Value = "product"
$wsInsured.Range("C27").Value = "email"
$wsInsured.Range("D27").Value = "brokerId"
$wsInsured.Range("E27").Value = "agentId"
$wsInsured.Range("F27").Value = "agencyOfficeId"
$wsInsured.Range("G27").Value = "functionality"
$wsInsured.Range("H27").Value = "clearanceText"

$wsInsured.Range("A28").Value = "Y"
$wsInsured.Range("B28").Value = "QA Program 5203"
$wsInsured.Range("C28").Value = "cfessler@profrisk.com"
$wsInsured.Range("D28").Value = 20217
$wsInsured.Range("E28").Value = 237
$wsInsured.Range("F28").Value = 8006
$wsInsured.Range("G28").Value = "submit"
$wsInsured.Range("H28").Value = "Test purpose"

# clear the leftover cells copied from the (wider) template that aren't part
# of this block: column I entirely, and C26:H26 (row 26 only has A26/B26)
$wsInsured.Range("I26:I28").Clear() | Out-Null
$wsInsured.Range("C26:H26").Clear() | Out-Null

# --- testCancelClearancesFunctionality (rows 31-33) ------------------------
$templateRange.Copy() | Out-Null
$wsInsured.Range("A31").PasteSpecial(-4122) | Out-Null

$wsInsured.Range("A31").Value = "testCancelClearancesFunctionality"

$wsInsured.Range("A32").Value = "runMode"
$wsInsured.Range("B32").Value = "product"
$wsInsured.Range("C32").Value = "email"
$wsInsured.Range("D32").Value = "brokerId"
$wsInsured.Range("E32").Value = "agentId"
$wsInsured.Range("F32").Value = "agencyOfficeId"
$wsInsured.Range("G32").Value = "functionality"
$wsInsured.Range("H32").Value = "clearanceText"

$wsInsured.Range("A33").Value = "Y"
$wsInsured.Range("B33").Value = "QA Program 5203"
$wsInsured.Range("C33").Value = "cfessler@profrisk.com"
$wsInsured.Range("D33").Value = 20217
$wsInsured.Range("E33").Value = 237
$wsInsured.Range("F33").Value = 8006
$wsInsured.Range("G33").Value = "cancel"
$wsInsured.Range("H33").Value = "Test purpose"

$wsInsured.Range("I31:I33").Clear() | Out-Null
$wsInsured.Range("C31:H31").Clear() | Out-Null

# --- testSubmissionClearanceComplete (rows 36-38) ---------------------------
$templateRange.Copy() | Out-Null
$wsInsured.Range("A36").PasteSpecial(-4122) | Out-Null

$wsInsured.Range("A36").Value = "testSubmissionClearanceComplete"

$wsInsured.Range("A37").Value = "runMode"
$wsInsured.Range("B37").Value = "product"
$wsInsured.Range("C37").Value = "applicantName"
$wsInsured.Range("D37").Value = "website"
$wsInsured.Range("E37").Value = "email"
$wsInsured.Range("F37").Value = "brokerId"
$wsInsured.Range("G37").Value = "agentId"
$wsInsured.Range("H37").Value = "agencyOfficeId"

$wsInsured.Range("A38").Value = "Y"
$wsInsured.Range("B38").Value = "QA Program 5203"
$wsInsured.Range("C38").Value = "venkatqa"
$wsInsured.Range("D38").Value = "venkatqa.com"
$wsInsured.Range("E38").Value = "cfessler@profrisk.com"
$wsInsured.Range("F38").Value = 20217
$wsInsured.Range("G38").Value = 237
$wsInsured.Range("H38").Value = 8006

$wsInsured.Range("I36:I38").Clear() | Out-Null
# row 36 only ever had column A filled in the source data
$wsInsured.Range("B36:H36").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 3. Make InsuredPageData the active sheet/tab, with the same selection the
#    workbook was saved with (D38)
# ---------------------------------------------------------------------------
$wsInsured.Activate() | Out-Null
$wsInsured.Range("D38").Select() | Out-Null
